$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the 12 now-unneeded blank timesheet rows (33-44) ---
# This shifts the old row 45 (totals row) up to row 33 and the old row 46
# (thick-top spacer row) up to row 34.
$ws.Rows("33:44").Delete()

# Row 32 (A32:F32) is now the last data row before the totals row, so it
# picks up the "closing" border treatment (no bottom border) that row 44
# used to have.
$ws.Range("E32:F32").Borders(9).LineStyle = -4142

# Re-assert the D-C shared formula over its new, shorter range so it keeps
# being stored as a single shared formula (F3:F32) instead of exploding
# into one formula per cell.
$ws.Range("F3:F32").Formula = "=D3-C3"

# --- Column width adjustments ---
$ws.Columns("A").ColumnWidth = 55.45
$ws.Columns("B").ColumnWidth = 11.1
$ws.Columns("F").ColumnWidth = 11.59

# --- Sheet view: drop the old scrolled/selected state, select the header row ---
$ws.Range("A1:F1").Select()

# --- Page setup: print landscape instead of portrait ---
$ws.PageSetup.Orientation = 2

Write-Host "done"
